# Auto-generated Excel COM-interop edit script
# Updates pricing/profit figures (columns H-N) for specific Leve rows
# across multiple job sheets, reflecting refreshed market-board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1839.6666
$ws.Range("I31").Value = 1445.25
$ws.Range("J31").Value = 4995
$ws.Range("K31").Value = 4335.75
$ws.Range("L31").Value = 14985
$ws.Range("M31").Value = -4105.75
$ws.Range("N31").Value = -15445

$ws.Range("H64").Value = 8006.8335
$ws.Range("I64").Value = 5017
$ws.Range("J64").Value = 10996.667
$ws.Range("K64").Value = 5017
$ws.Range("L64").Value = 10996.667
$ws.Range("M64").Value = -4769
$ws.Range("N64").Value = -11492.667

$ws.Range("H67").Value = 8006.8335
$ws.Range("I67").Value = 5017
$ws.Range("J67").Value = 10996.667
$ws.Range("K67").Value = 5017
$ws.Range("L67").Value = 10996.667
$ws.Range("M67").Value = -4159
$ws.Range("N67").Value = -12712.667

$ws.Range("H88").Value = 8699
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 8699
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 8699
$ws.Range("N88").Value = -9511

$ws.Range("H91").Value = 8699
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 8699
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 8699
$ws.Range("N91").Value = -11507

$ws.Range("H101").Value = 1094.8462
$ws.Range("I101").Value = 1065.4445
$ws.Range("J101").Value = 1161
$ws.Range("K101").Value = 3196.3335
$ws.Range("L101").Value = 3483
$ws.Range("M101").Value = -1574.3335
$ws.Range("N101").Value = -6727

$ws.Range("H112").Value = 1883.7273
$ws.Range("I112").Value = 1640
$ws.Range("J112").Value = 1955.4117
$ws.Range("K112").Value = 4920
$ws.Range("L112").Value = 5866.2351
$ws.Range("M112").Value = -3812
$ws.Range("N112").Value = -8082.2351

$ws.Range("H113").Value = 2502.25
$ws.Range("I113").Value = 1670
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 1670
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = 1584
$ws.Range("N113").Value = -11507

$ws.Range("H129").Value = 3025.1875
$ws.Range("I129").Value = 1402.3636
$ws.Range("J129").Value = 6595.4
$ws.Range("K129").Value = 4207.0908
$ws.Range("L129").Value = 19786.2
$ws.Range("M129").Value = 792.9092000000001
$ws.Range("N129").Value = -29786.2

$ws.Range("H138").Value = 1936.5178
$ws.Range("I138").Value = 1013.8485
$ws.Range("J138").Value = 3260.348
$ws.Range("K138").Value = 3041.5455
$ws.Range("L138").Value = 9781.044
$ws.Range("M138").Value = 2098.4545
$ws.Range("N138").Value = -20061.044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 44393.406
$ws.Range("I74").Value = 29629.688
$ws.Range("J74").Value = 91637.3
$ws.Range("K74").Value = 29629.688
$ws.Range("L74").Value = 91637.3
$ws.Range("M74").Value = -28755.688
$ws.Range("N74").Value = -93385.3

$ws.Range("H77").Value = 44393.406
$ws.Range("I77").Value = 29629.688
$ws.Range("J77").Value = 91637.3
$ws.Range("K77").Value = 148148.44
$ws.Range("L77").Value = 458186.5
$ws.Range("M77").Value = -143780.44
$ws.Range("N77").Value = -466922.5

$ws.Range("H132").Value = 1972.4722
$ws.Range("I132").Value = 1894
$ws.Range("J132").Value = 2297.5715
$ws.Range("K132").Value = 5682
$ws.Range("L132").Value = 6892.7145
$ws.Range("M132").Value = -3152
$ws.Range("N132").Value = -11952.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("N22").Value = 0

$ws.Range("H29").Value = 823
$ws.Range("I29").Value = 823
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 823
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -534

$ws.Range("H36").Value = 14995
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 14995
$ws.Range("K36").Value = 0
$ws.Range("L36").ClearContents()
$ws.Range("M36").Value = 14995
$ws.Range("N36").Value = -16063

$ws.Range("H99").Value = 2265.625
$ws.Range("I99").Value = 2346.7144
$ws.Range("J99").Value = 1698
$ws.Range("K99").Value = 2346.7144
$ws.Range("L99").Value = 1698
$ws.Range("M99").Value = -848.7143999999998
$ws.Range("N99").Value = -4694

$ws.Range("H134").Value = 1996.8654
$ws.Range("I134").Value = 1605.75
$ws.Range("J134").Value = 3300.5833
$ws.Range("K134").Value = 4817.25
$ws.Range("L134").Value = 9901.749899999999
$ws.Range("M134").Value = -2282.25
$ws.Range("N134").Value = -14971.7499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 281.36365
$ws.Range("I7").Value = 204.75
$ws.Range("J7").Value = 325.14285
$ws.Range("K7").Value = 204.75
$ws.Range("L7").Value = 325.14285
$ws.Range("M7").Value = -91.75
$ws.Range("N7").Value = -551.14285

$ws.Range("H22").Value = 1109.2
$ws.Range("I22").Value = 490.85715
$ws.Range("J22").Value = 1650.25
$ws.Range("K22").Value = 490.85715
$ws.Range("L22").Value = 1650.25
$ws.Range("M22").Value = -140.85715
$ws.Range("N22").Value = -2350.25

$ws.Range("H51").Value = 46999.332
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 46999.332
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 46999.332
$ws.Range("N51").Value = -48471.332

$ws.Range("H58").Value = 1594.3846
$ws.Range("I58").Value = 1140.5625
$ws.Range("J58").Value = 3669
$ws.Range("K58").Value = 1140.5625
$ws.Range("L58").Value = 3669
$ws.Range("M58").Value = -937.5625
$ws.Range("N58").Value = -4075

$ws.Range("H61").Value = 46999.332
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 46999.332
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 46999.332
$ws.Range("N61").Value = -47695.332

$ws.Range("H132").Value = 3657
$ws.Range("I132").Value = 3131.5
$ws.Range("J132").Value = 5338.6
$ws.Range("K132").Value = 9394.5
$ws.Range("L132").Value = 16015.8
$ws.Range("M132").Value = -6864.5
$ws.Range("N132").Value = -21075.8

$ws.Range("H134").Value = 2140.5
$ws.Range("I134").Value = 2180.8254
$ws.Range("J134").Value = 1293.6666
$ws.Range("K134").Value = 6542.476200000001
$ws.Range("L134").Value = 3880.9998
$ws.Range("M134").Value = -4007.476200000001
$ws.Range("N134").Value = -8950.9998

$ws.Range("H136").Value = 1594.3846
$ws.Range("I136").Value = 1140.5625
$ws.Range("J136").Value = 3669
$ws.Range("K136").Value = 3421.6875
$ws.Range("L136").Value = 11007
$ws.Range("M136").Value = -871.6875
$ws.Range("N136").Value = -16107

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1815.3334
$ws.Range("I132").Value = 1717.9231
$ws.Range("J132").Value = 1973.625
$ws.Range("K132").Value = 15461.3079
$ws.Range("L132").Value = 17762.625
$ws.Range("M132").Value = -12931.3079
$ws.Range("N132").Value = -22822.625

$ws.Range("H137").Value = 1882.7778
$ws.Range("I137").Value = 1173.5625
$ws.Range("J137").Value = 2914.3635
$ws.Range("K137").Value = 3520.6875
$ws.Range("L137").Value = 8743.0905
$ws.Range("M137").Value = 1579.3125
$ws.Range("N137").Value = -18943.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 969.73334
$ws.Range("I113").Value = 892.55554
$ws.Range("J113").Value = 1085.5
$ws.Range("K113").Value = 892.55554
$ws.Range("L113").Value = 1085.5
$ws.Range("M113").Value = 1277.44446
$ws.Range("N113").Value = -5425.5

$ws.Range("H122").Value = 106120.48
$ws.Range("I122").Value = 167402.92
$ws.Range("J122").Value = 1064.8572
$ws.Range("K122").Value = 502208.76
$ws.Range("L122").Value = 3194.5716
$ws.Range("M122").Value = -499758.76
$ws.Range("N122").Value = -8094.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22802.408
$ws.Range("I7").Value = 26248.588
$ws.Range("J7").Value = 11085.4
$ws.Range("K7").Value = 26248.588
$ws.Range("L7").Value = 11085.4
$ws.Range("M7").Value = -26136.588
$ws.Range("N7").Value = -11309.4

$ws.Range("H21").Value = 8157.375
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 8157.375
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 8157.375
$ws.Range("N21").Value = -8505.375

$ws.Range("H40").Value = 25136.945
$ws.Range("I40").Value = 27967.062
$ws.Range("J40").Value = 2496
$ws.Range("K40").Value = 27967.062
$ws.Range("L40").Value = 2496
$ws.Range("M40").Value = -27831.062
$ws.Range("N40").Value = -2768

$ws.Range("H82").Value = 1554.6818
$ws.Range("I82").Value = 1481.8
$ws.Range("J82").Value = 1710.8572
$ws.Range("K82").Value = 1481.8
$ws.Range("L82").Value = 1710.8572
$ws.Range("M82").Value = -1120.8
$ws.Range("N82").Value = -2432.8572

$ws.Range("H85").Value = 1554.6818
$ws.Range("I85").Value = 1481.8
$ws.Range("J85").Value = 1710.8572
$ws.Range("K85").Value = 1481.8
$ws.Range("L85").Value = 1710.8572
$ws.Range("M85").Value = -233.8
$ws.Range("N85").Value = -4206.8572

$ws.Range("H126").Value = 22802.408
$ws.Range("I126").Value = 26248.588
$ws.Range("J126").Value = 11085.4
$ws.Range("K126").Value = 78745.764
$ws.Range("L126").Value = 33256.2
$ws.Range("M126").Value = -76275.764
$ws.Range("N126").Value = -38196.2

$ws.Range("H132").Value = 2642.1777
$ws.Range("I132").Value = 2235.2163
$ws.Range("J132").Value = 4524.375
$ws.Range("K132").Value = 6705.6489
$ws.Range("L132").Value = 13573.125
$ws.Range("M132").Value = -4175.6489
$ws.Range("N132").Value = -18633.125

$ws.Range("H136").Value = 26812.838
$ws.Range("I136").Value = 2139.3928
$ws.Range("J136").Value = 103574.664
$ws.Range("K136").Value = 6418.178400000001
$ws.Range("L136").Value = 310723.992
$ws.Range("M136").Value = -3868.178400000001
$ws.Range("N136").Value = -315823.992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 17499.5
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 24999
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 24999
$ws.Range("M20").Value = -9760
$ws.Range("N20").Value = -25479

$ws.Range("H62").Value = 17829.846
$ws.Range("I62").Value = 20448.5
$ws.Range("J62").Value = 16666
$ws.Range("K62").Value = 20448.5
$ws.Range("L62").Value = 16666
$ws.Range("M62").Value = -19824.5
$ws.Range("N62").Value = -17914

$ws.Range("H65").Value = 17829.846
$ws.Range("I65").Value = 20448.5
$ws.Range("J65").Value = 16666
$ws.Range("K65").Value = 102242.5
$ws.Range("L65").Value = 83330
$ws.Range("M65").Value = -99122.5
$ws.Range("N65").Value = -89570

$ws.Range("H122").Value = 2836.0588
$ws.Range("I122").Value = 2663.0356
$ws.Range("J122").Value = 3643.5
$ws.Range("K122").Value = 7989.1068
$ws.Range("L122").Value = 10930.5
$ws.Range("M122").Value = -5539.1068
$ws.Range("N122").Value = -15830.5

$ws.Range("H126").Value = 83338760
$ws.Range("I126").Value = 100005310
$ws.Range("J126").Value = 5999.5
$ws.Range("K126").Value = 300015930
$ws.Range("L126").Value = 17998.5
$ws.Range("M126").Value = -300013460
$ws.Range("N126").Value = -22938.5
